$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 11.14494766666667
$ws.Range("H2").Value = 33.434843
$ws.Range("I2").Value = 0.1279818847384872
$ws.Range("J2").Value = 0.1279818847384872
$ws.Range("M2").Value = 0.2819746666666667
$ws.Range("N2").Value = 0.845924
$ws.Range("O2").Value = 0.02777466224158379
$ws.Range("P2").Value = 0.02777466224158379
$ws.Range("Q2").Value = 3.142592903325778
$ws.Range("R2").Value = 28.283336129932
$ws.Range("S2").Value = 0.00355465362165279
$ws.Range("T2").Value = 0.00355465362165279
$ws.Range("G3").Value = 11.14494766666667
$ws.Range("H3").Value = 33.434843
$ws.Range("I3").Value = 0.1279818847384872
$ws.Range("J3").Value = 0.1279818847384872
$ws.Range("O3").Value = 0.8500493487799721
$ws.Range("P3").Value = 0.8500493487799721
$ws.Range("Q3").Value = 96.17971328389811
$ws.Range("R3").Value = 865.617419555083
$ws.Range("S3").Value = 0.1087909177775845
$ws.Range("T3").Value = 0.1087909177775845
$ws.Range("G4").Value = 11.14494766666667
$ws.Range("H4").Value = 33.434843
$ws.Range("I4").Value = 0.1279818847384872
$ws.Range("J4").Value = 0.1279818847384872
$ws.Range("M4").Value = 1.240358333333333
$ws.Range("N4").Value = 3.721075
$ws.Range("O4").Value = 0.1221759889784441
$ws.Range("P4").Value = 0.1221759889784441
$ws.Range("Q4").Value = 13.82372871291389
$ws.Range("R4").Value = 124.413558416225
$ws.Range("S4").Value = 0.01563631333924993
$ws.Range("T4").Value = 0.01563631333924993
$ws.Range("I5").Value = 0.5307607770439682
$ws.Range("J5").Value = 0.5307607770439681
$ws.Range("M5").Value = 0.2819746666666667
$ws.Range("N5").Value = 0.845924
$ws.Range("O5").Value = 0.02777466224158379
$ws.Range("P5").Value = 0.02777466224158379
$ws.Range("Q5").Value = 13.03282143961466
$ws.Range("R5").Value = 117.295392956532
$ws.Range("S5").Value = 0.01474170131347677
$ws.Range("T5").Value = 0.01474170131347677
$ws.Range("I6").Value = 0.5307607770439682
$ws.Range("J6").Value = 0.5307607770439681
$ws.Range("O6").Value = 0.8500493487799721
$ws.Range("P6").Value = 0.8500493487799721
$ws.Range("S6").Value = 0.4511728528841771
$ws.Range("T6").Value = 0.451172852884177
$ws.Range("I7").Value = 0.5307607770439682
$ws.Range("J7").Value = 0.5307607770439681
$ws.Range("M7").Value = 1.240358333333333
$ws.Range("N7").Value = 3.721075
$ws.Range("O7").Value = 0.1221759889784441
$ws.Range("P7").Value = 0.1221759889784441
$ws.Range("Q7").Value = 57.32915254610833
$ws.Range("R7").Value = 515.9623729149749
$ws.Range("S7").Value = 0.06484622284631431
$ws.Range("T7").Value = 0.06484622284631429
$ws.Range("G8").Value = 29.71744933333333
$ws.Range("H8").Value = 89.152348
$ws.Range("I8").Value = 0.3412573382175446
$ws.Range("J8").Value = 0.3412573382175446
$ws.Range("M8").Value = 0.2819746666666667
$ws.Range("N8").Value = 0.845924
$ws.Range("O8").Value = 0.02777466224158379
$ws.Range("P8").Value = 0.02777466224158379
$ws.Range("Q8").Value = 8.379567869950222
$ws.Range("R8").Value = 75.41611082955201
$ws.Range("S8").Value = 0.009478307306454223
$ws.Range("T8").Value = 0.009478307306454223
$ws.Range("G9").Value = 29.71744933333333
$ws.Range("H9").Value = 89.152348
$ws.Range("I9").Value = 0.3412573382175446
$ws.Range("J9").Value = 0.3412573382175446
$ws.Range("O9").Value = 0.8500493487799721
$ws.Range("P9").Value = 0.8500493487799721
$ws.Range("Q9").Value = 256.4584277912209
$ws.Range("R9").Value = 2308.125850120988
$ws.Range("S9").Value = 0.2900855781182105
$ws.Range("T9").Value = 0.2900855781182105
$ws.Range("G10").Value = 29.71744933333333
$ws.Range("H10").Value = 89.152348
$ws.Range("I10").Value = 0.3412573382175446
$ws.Range("J10").Value = 0.3412573382175446
$ws.Range("M10").Value = 1.240358333333333
$ws.Range("N10").Value = 3.721075
$ws.Range("O10").Value = 0.1221759889784441
$ws.Range("P10").Value = 0.1221759889784441
$ws.Range("Q10").Value = 36.86028592601111
$ws.Range("R10").Value = 331.7425733341
$ws.Range("S10").Value = 0.04169345279287991
$ws.Range("T10").Value = 0.04169345279287991
